{"js": "// Insert \"22\" before \"students, Layout 2 was chosen\" in the \"Student\n// Preferences\" paragraph, so the sentence reads: \"...gathered from the 22\n// students, Layout 2 was chosen...\" (author's commit message: \"added\n// amount of student tested\").\n//\n// Note: \"students, Layout 2 was chosen\" is unique in the document (a very\n// similar, but distinct, sentence about the \"General Usability Test\"\n// appears later and must stay untouched), so searching for it directly is\n// safe and precise.\n\nconst body = context.document.body;\n\nconst searchResults = body.search(\"students, Layout 2 was chosen\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Anchor text \"students, Layout 2 was chosen\" not found.');\n}\n\n// Insert \"22 \" immediately before \"students\" -> \"...the 22 students, ...\".\nsearchResults.items[0].insertText(\"22 \", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Insert \"22\" before \"students, Layout 2 was chosen\" in the \"Student\n# Preferences\" paragraph, so the sentence reads: \"...gathered from the 22\n# students, Layout 2 was chosen...\" (author's commit message: \"added\n# amount of student tested\").\n\n$d = $word.ActiveDocument\n\n# Use a narrow, unique anchor so we don't touch the other, near-identical\n# \"students, Layout 2\" sentence later in the document.\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"students, Layout 2 was chosen\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif ($found -and $find.Found) {\n    # After Execute(), $rng has been collapsed/moved to the matched text\n    # (\"students, Layout 2 was chosen\"). Insert \"22 \" right before it so the\n    # sentence becomes \"...the 22 students, Layout 2 was chosen...\".\n    $rng.InsertBefore(\"22 \")\n}\n"}
